$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.987569333333333
$ws.Range("H2").Value = 5.962707999999999
$ws.Range("I2").Value = 0.1241595834663642
$ws.Range("J2").Value = 0.1241595834663642
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 105.1628969852209
$ws.Range("R2").Value = 946.4660728669878
$ws.Range("S2").Value = 0.0516679605868599
$ws.Range("T2").Value = 0.05166796058685991

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.987569333333333
$ws.Range("H3").Value = 5.962707999999999
$ws.Range("I3").Value = 0.1241595834663642
$ws.Range("J3").Value = 0.1241595834663642
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 93.93264263487775
$ws.Range("R3").Value = 845.3937837138998
$ws.Range("S3").Value = 0.04615038399104317
$ws.Range("T3").Value = 0.04615038399104318

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.987569333333333
$ws.Range("H4").Value = 5.962707999999999
$ws.Range("I4").Value = 0.1241595834663642
$ws.Range("J4").Value = 0.1241595834663642
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 53.61390231444177
$ws.Range("R4").Value = 482.5251208299759
$ws.Range("S4").Value = 0.02634123888846115
$ws.Range("T4").Value = 0.02634123888846116

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.481595333333334
$ws.Range("H5").Value = 22.444786
$ws.Range("I5").Value = 0.4673606825542495
$ws.Range("J5").Value = 0.4673606825542495
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 395.8534809977829
$ws.Range("R5").Value = 3562.681328980046
$ws.Range("S5").Value = 0.1944881953683637
$ws.Range("T5").Value = 0.1944881953683637

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.481595333333334
$ws.Range("H6").Value = 22.444786
$ws.Range("I6").Value = 0.4673606825542495
$ws.Range("J6").Value = 0.4673606825542495
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 353.5806318797278
$ws.Range("R6").Value = 3182.22568691755
$ws.Range("S6").Value = 0.1737189700546782
$ws.Range("T6").Value = 0.1737189700546782

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.481595333333334
$ws.Range("H7").Value = 22.444786
$ws.Range("I7").Value = 0.4673606825542495
$ws.Range("J7").Value = 0.4673606825542495
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 201.8130963435658
$ws.Range("R7").Value = 1816.317867092092
$ws.Range("S7").Value = 0.0991535171312076
$ws.Range("T7").Value = 0.09915351713120761

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.539018333333334
$ws.Range("H8").Value = 19.617055
$ws.Range("I8").Value = 0.4084797339793862
$ws.Range("J8").Value = 0.4084797339793863
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 345.9814456985673
$ws.Range("R8").Value = 3113.833011287105
$ws.Range("S8").Value = 0.1699853866012327
$ws.Range("T8").Value = 0.1699853866012327

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.539018333333334
$ws.Range("H9").Value = 19.617055
$ws.Range("I9").Value = 0.4084797339793862
$ws.Range("J9").Value = 0.4084797339793863
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 309.0343878760694
$ws.Range("R9").Value = 2781.309490884625
$ws.Range("S9").Value = 0.1518327949353571
$ws.Range("T9").Value = 0.1518327949353572

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.539018333333334
$ws.Range("H10").Value = 19.617055
$ws.Range("I10").Value = 0.4084797339793862
$ws.Range("J10").Value = 0.4084797339793863
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 176.3874518871345
$ws.Range("R10").Value = 1587.48706698421
$ws.Range("S10").Value = 0.08666155244279636
$ws.Range("T10").Value = 0.08666155244279637

